# Weekly update: add two new "Arándano (blue)" price records for
# Vega Central Mapocho de Santiago, inserted right above the existing
# row 103 (pushing the rest of the table down by two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 103; this shifts old rows 103:113
# down to 105:115 (and carries the date-format style from row 103
# onto the new rows, matching Excel's default insert behaviour).
$ws.Rows.Item(103).Insert()
$ws.Rows.Item(103).Insert()

# New row 103
$ws.Range("A103").Value = 9
$ws.Range("B103").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C103").Value = "Metropolitana"
$ws.Range("D103").Value = 44504
$ws.Range("E103").Value = 13
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100101
$ws.Range("H103").Value = "Berries"
$ws.Range("I103").Value = 100101001
$ws.Range("J103").Value = "Arándano (blue)"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 350
$ws.Range("N103").Value = 14000
$ws.Range("O103").Value = 14000
$ws.Range("P103").Value = 14000
$ws.Range("Q103").Value = "$/bandeja 2 kilos"
$ws.Range("R103").Value = "Provincia de Curicó"
$ws.Range("S103").Value = 7000
$ws.Range("T103").Value = 2

# New row 104
$ws.Range("A104").Value = 9
$ws.Range("B104").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C104").Value = "Metropolitana"
$ws.Range("D104").Value = 44504
$ws.Range("E104").Value = 13
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100101
$ws.Range("H104").Value = "Berries"
$ws.Range("I104").Value = 100101001
$ws.Range("J104").Value = "Arándano (blue)"
$ws.Range("K104").Value = "Sin especificar"
$ws.Range("L104").Value = "Segunda"
$ws.Range("M104").Value = 310
$ws.Range("N104").Value = 12000
$ws.Range("O104").Value = 12000
$ws.Range("P104").Value = 12000
$ws.Range("Q104").Value = "$/bandeja 2 kilos"
$ws.Range("R104").Value = "Provincia de Curicó"
$ws.Range("S104").Value = 6000
$ws.Range("T104").Value = 2
